$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Insert a new row right after the header row (i.e. before the
#     existing "9/28/2022" row), for the 10/25/2022 death-data update. ---
$firstDataRow = $t.Rows.Item(2)
$newRow = $t.Rows.Add($firstDataRow)

$t.Cell(2,1).Range.InsertXML("<w:p $wns><w:r><w:t>10/25/2022</w:t></w:r></w:p>")

$t.Cell(2,2).Range.InsertXML("<w:p $wns><w:r><w:t>All CCB + ad-hoc death datasets</w:t></w:r></w:p>")

$reasonXml = "<w:p $wns>" +
  "<w:r><w:t>Received regularly scheduled death data from CHSI on 2</w:t></w:r>" +
  "<w:r><w:rPr><w:vertAlign w:val=""superscript""/></w:rPr><w:t>nd</w:t></w:r>" +
  "<w:r><w:t xml:space=""preserve""> week of October 2022. This round of death data received consists of the final 2021 CCDF file, and preliminary year-to-date death data.</w:t></w:r>" +
  "</w:p>"
$t.Cell(2,3).Range.InsertXML($reasonXml)

$t.Cell(2,4).Range.InsertXML("<w:p $wns><w:r><w:t>0.Secure.Data/ Archived Data/ archiveDat/ 20221025</w:t></w:r></w:p>")

# --- In the (now third) row, "9/28/2022", split the "Datasets updated"
#     cell text into three runs, inserting "death " in the middle. ---
$datasetsXml = "<w:p $wns>" +
  "<w:r><w:t xml:space=""preserve"">All CCB + ad-hoc </w:t></w:r>" +
  "<w:r><w:t xml:space=""preserve"">death </w:t></w:r>" +
  "<w:r><w:t>datasets</w:t></w:r>" +
  "</w:p>"
$t.Cell(3,2).Range.InsertXML($datasetsXml)
